$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title of the left (C) table -----------------------------------------
$ws.Range("A1").Value = "C, primes until 1750, upper limit: 2^1000000"

# --- Updated execution times for the C run (column B) --------------------
$ws.Range("B3").Value = 195.26
$ws.Range("B4").Value = 244.2
$ws.Range("B5").Value = 291.37
$ws.Range("B6").Value = 369.59
$ws.Range("B7").Value = 364.32
$ws.Range("B8").Value = 376.23
$ws.Range("B9").Value = 395.68
$ws.Range("B10").Value = 388.31

# --- New right-hand (C++) table header ------------------------------------
$ws.Range("F1:I1").Merge()
$ws.Range("F1").Value = "C++, primes until 1750, upper limit: 2^1000000"
$ws.Range("F1:I1").HorizontalAlignment = -4108

$ws.Range("F2").Value = "Processes"
$ws.Range("G2").Value = "Execution time (s)"
$ws.Range("H2").Value = "Speedup"
$ws.Range("I2").Value = "Efficiency"

# --- New right-hand (C++) table data --------------------------------------
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 332.596

$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 215.139
$ws.Range("H4").Formula = "=G3/G4"
$ws.Range("I4").Formula = "=H4/F4"

$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 201.66
$ws.Range("H5").Formula = "=G3/G5"
$ws.Range("I5").Formula = "=H5/F5"

$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 189.026
$ws.Range("H6").Formula = "=G3/G6"
$ws.Range("I6").Formula = "=H6/F6"

$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 186.32
$ws.Range("H7").Formula = "=G3/G7"
$ws.Range("I7").Formula = "=H7/F7"

$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 184.761
$ws.Range("H8").Formula = "=G3/G8"
$ws.Range("I8").Formula = "=H8/F8"

$ws.Range("F9").Value = 7
$ws.Range("G9").Value = 188.008
$ws.Range("H9").Formula = "=G3/G9"
$ws.Range("I9").Formula = "=H9/F9"

$ws.Range("F10").Value = 8
$ws.Range("G10").Value = 183.372
$ws.Range("H10").Formula = "=G3/G10"
$ws.Range("I10").Formula = "=H10/F10"

# --- Column widths ----------------------------------------------------------
# C and D widen to fit the new, longer titles
$ws.Columns.Item(3).ColumnWidth = 11.17
$ws.Columns.Item(4).ColumnWidth = 11.17

# F/G mirror A/B, H/I mirror C/D
$ws.Columns.Item(6).ColumnWidth = 8.83
$ws.Columns.Item(7).ColumnWidth = 16.5
$ws.Columns.Item(8).ColumnWidth = 11.17
$ws.Columns.Item(9).ColumnWidth = 11.17

# --- Selection, matching the author's final cursor position ---------------
$ws.Range("M11").Select()
